# "desglosse del 15 al 29 de feb"
# Add a third worksheet "Desglose de ingresos 15-29 Feb" (a copy of the
# first sheet's layout/formulas/styles) with its own figures, make it the
# active sheet/tab, and update the now-inactive first sheet's selection.

$wb = $excel.ActiveWorkbook

# --- 1. First sheet ("15-31 Ener") loses the active-tab flag and its
#        selection moves to B6 (it is no longer the selected tab). ---
$ws1 = $wb.Worksheets.Item(1)
[void]$ws1.Activate()
[void]$ws1.Range("B6").Select()

# --- 2. Duplicate sheet1 to the end of the workbook: this clones its
#        columns/styles/formulas (same style ids, same SUM formulas) so
#        the new sheet starts out identical, ready to receive new figures.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
[void]$ws1.Copy($null, $lastSheet)
$ws3 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3.Name = "Desglose de ingresos 15-29 Feb"

# --- 3. Overwrite the figures for the 15-29 Feb breakdown. ---
$ws3.Range("A2").Value = 149.47
$ws3.Range("C2").Value = 7846.06
$ws3.Range("E2").Value = 7846.06

$ws3.Range("A3").Value = 2145
$ws3.Range("C3").Value = 1202.69
$ws3.Range("E3").Value = 1202.69

$ws3.Range("A4").Value = 71.15

$ws3.Range("A5").Value = 53.99

$ws3.Range("A6").Value = 53.99

$ws3.Range("A9").Value = 2116.01
$ws3.Range("C9").Value = 7846.06
$ws3.Range("E9").Value = 7846.06

# --- 4. The new sheet becomes the active tab/selection. ---
[void]$ws3.Range("I8").Select()
